# Fall 2021 schedule update: replace plain lecture-topic text in column C
# with markdown-style links pointing at the online textbook pages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value  = "[What is Cognition?](https://www.crumplab.com/cognition/textbook/what-is-cognition.html)"
$ws.Range("C4").Value  = "[Mental Imagery](https://www.crumplab.com/cognition/textbook/mental-imagery.html)"
$ws.Range("C7").Value  = "[Eugenics and Psychology](https://www.crumplab.com/cognition/textbook/eugenics-and-psychology.html)"
$ws.Range("C9").Value  = "[Intelligence Testing](https://www.crumplab.com/cognition/textbook/intelligence-testing.html)"
$ws.Range("C10").Value = "[Associations](https://www.crumplab.com/cognition/textbook/associations.html)"
$ws.Range("C11").Value = "[Associations](https://www.crumplab.com/cognition/textbook/associations.html)"
$ws.Range("C12").Value = "[Behaviorism](https://www.crumplab.com/cognition/textbook/behaviorism.html)"
$ws.Range("C13").Value = "[Information Processing](https://www.crumplab.com/cognition/textbook/information-processing.html)"
$ws.Range("C14").Value = "[Information Processing](https://www.crumplab.com/cognition/textbook/information-processing.html)"
$ws.Range("C17").Value = "[Memory I](https://www.crumplab.com/cognition/textbook/memory-i.html)"
$ws.Range("C18").Value = "[Memory I](https://www.crumplab.com/cognition/textbook/memory-i.html)"
$ws.Range("C19").Value = "[Memory II](https://www.crumplab.com/cognition/textbook/memory-ii.html)"
$ws.Range("C20").Value = "[Memory II](https://www.crumplab.com/cognition/textbook/memory-ii.html)"

# Reflect the view state recorded in the saved workbook (window geometry
# and current selection) as captured when the author saved the file.
$ws.Range("C21").Select()
$ws.Application.ActiveWindow.ScrollRow = 4

$excel.Width = 17040
$excel.Height = 20540
$excel.Left = 7040
$excel.Top = 460
